$d = $word.ActiveDocument

# The two Pearson logo images (one in the "first page" footer, one in the
# "default" footer) are renamed from image1.png -> image2.png, and the
# BTEC logo image in the "first page" header is renamed from image2.jpg ->
# image1.jpg. Walk every header/footer of every section and rename the
# inline picture(s) found there based on which logo they are (identified
# by their AlternativeText / descr, which the edit leaves untouched).
foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
